# Weekly update: a new day's price observation for Pomelo (Start Ruby,
# Primera) at "Vega Modelo de Temuco" is inserted as a new row 346,
# pushing all subsequent rows (old 346..443) down by one (new 347..444).
#
# This mirrors the canonical OOXML diff exactly: every row from 347
# onward ends up holding the data that used to live one row above it,
# the worksheet's used range grows from A1:T443 to A1:T444, and the
# brand-new row 346 carries the newly reported observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 346..443 down to 347..444, creating a blank row 346.
$ws.Rows.Item(346).Insert()

# Populate the newly inserted row 346 with the new observation.
$ws.Cells.Item(346, 1).Value  = 10
$ws.Cells.Item(346, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(346, 3).Value  = "La Araucanía"
$ws.Cells.Item(346, 4).Value  = 45093
$ws.Cells.Item(346, 5).Value  = 9
$ws.Cells.Item(346, 6).Value  = "Fruta"
$ws.Cells.Item(346, 7).Value  = 100102
$ws.Cells.Item(346, 8).Value  = "Cítricos"
$ws.Cells.Item(346, 9).Value  = 100102006
$ws.Cells.Item(346, 10).Value = "Pomelo"
$ws.Cells.Item(346, 11).Value = "Start Ruby"
$ws.Cells.Item(346, 12).Value = "Primera"
$ws.Cells.Item(346, 13).Value = 170
$ws.Cells.Item(346, 14).Value = 14000
$ws.Cells.Item(346, 15).Value = 15000
$ws.Cells.Item(346, 16).Value = 14471
$ws.Cells.Item(346, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(346, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(346, 19).Value = 965
$ws.Cells.Item(346, 20).Value = 15
